$wb = $excel.ActiveWorkbook

# ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 5320.6665
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

# ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 5320.6665
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

# ALC row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 874.125
$ws.Range("I86").Value = 874.125
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 874.125
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 248.875
$ws.Range("N86").ClearContents()

# ALC row 88
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 670978.8
$ws.Range("J88").Value = 913305
$ws.Range("L88").Value = 913305
$ws.Range("N88").Value = -914117

# ALC row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 874.125
$ws.Range("I89").Value = 874.125
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 4370.625
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = 1245.375
$ws.Range("N89").ClearContents()

# ALC row 91
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 670978.8
$ws.Range("J91").Value = 913305
$ws.Range("L91").Value = 913305
$ws.Range("N91").Value = -916113

# ALC row 107
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 2989.2307
$ws.Range("I107").Value = 1735.7
$ws.Range("J107").Value = 7167.6665
$ws.Range("K107").Value = 1735.7
$ws.Range("L107").Value = 7167.6665
$ws.Range("M107").Value = 184.3
$ws.Range("N107").Value = -11007.6665

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2054.25
$ws.Range("I138").Value = 1354.8462
$ws.Range("J138").Value = 2880.818
$ws.Range("K138").Value = 4064.5386
$ws.Range("L138").Value = 8642.454000000002
$ws.Range("M138").Value = 1075.4614
$ws.Range("N138").Value = -18922.454

# ARM row 23
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 238000.2
$ws.Range("I23").Value = 22500
$ws.Range("K23").Value = 22500
$ws.Range("M23").Value = -22241

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3763.2703
$ws.Range("I32").Value = 2796.4849
$ws.Range("J32").Value = 11739.25
$ws.Range("K32").Value = 2796.4849
$ws.Range("L32").Value = 11739.25
$ws.Range("M32").Value = -2509.4849
$ws.Range("N32").Value = -12313.25

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1313.3334
$ws.Range("I110").Value = 1335.7142
$ws.Range("K110").Value = 1335.7142
$ws.Range("M110").Value = 709.2858000000001

# BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 715.2
$ws.Range("I20").Value = 615.38464
$ws.Range("J20").Value = 1364
$ws.Range("K20").Value = 615.38464
$ws.Range("L20").Value = 1364
$ws.Range("M20").Value = -368.38464
$ws.Range("N20").Value = -1858

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 5348.7
$ws.Range("I105").Value = 5348.7
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 5348.7
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -3601.7
$ws.Range("N105").ClearContents()

# CRP row 4
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 5817.8184
$ws.Range("I4").Value = 1000
$ws.Range("K4").Value = 1000
$ws.Range("M4").Value = -888

# CRP row 6
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 2143028.5
$ws.Range("I6").Value = 2143028.5
$ws.Range("K6").Value = 2143028.5
$ws.Range("M6").Value = -2142915.5

# CRP row 19
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 1539699.5
$ws.Range("I19").Value = 1667924.5
$ws.Range("J19").Value = 999
$ws.Range("K19").Value = 1667924.5
$ws.Range("L19").Value = 999
$ws.Range("M19").Value = -1667754.5
$ws.Range("N19").Value = -1339

# CRP row 24
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H24").Value = 1539699.5
$ws.Range("I24").Value = 1667924.5
$ws.Range("J24").Value = 999
$ws.Range("K24").Value = 1667924.5
$ws.Range("L24").Value = 999
$ws.Range("M24").Value = -1667754.5
$ws.Range("N24").Value = -1339

# CRP row 74
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 56999.855
$ws.Range("J74").Value = 56999.855
$ws.Range("L74").Value = 56999.855
$ws.Range("N74").Value = -58747.855

# CRP row 77
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 56999.855
$ws.Range("J77").Value = 56999.855
$ws.Range("L77").Value = 170999.565
$ws.Range("N77").Value = -179735.565

# CRP row 88
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 10945
$ws.Range("J88").Value = 10945
$ws.Range("L88").Value = 10945
$ws.Range("N88").Value = -11757

# CRP row 91
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H91").Value = 10945
$ws.Range("J91").Value = 10945
$ws.Range("L91").Value = 10945
$ws.Range("N91").Value = -13753

# CUL row 34
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2734.0908
$ws.Range("I34").Value = 190
$ws.Range("J34").Value = 3299.4443
$ws.Range("K34").Value = 570
$ws.Range("L34").Value = 9898.332900000001
$ws.Range("M34").Value = -486
$ws.Range("N34").Value = -10066.3329

# CUL row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1078.5883
$ws.Range("I107").Value = 1114.7778
$ws.Range("J107").Value = 1037.875
$ws.Range("K107").Value = 3344.3334
$ws.Range("L107").Value = 3113.625
$ws.Range("M107").Value = -1424.3334
$ws.Range("N107").Value = -6953.625

# CUL row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 927.9474
$ws.Range("I140").Value = 927.9474
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 2783.8422
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = 2396.1578
$ws.Range("N140").ClearContents()

# GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 509.47058
$ws.Range("I97").Value = 511
$ws.Range("K97").Value = 511
$ws.Range("M97").Value = -15

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4024.3333
$ws.Range("I22").Value = 3399.25
$ws.Range("J22").Value = 4171.4116
$ws.Range("K22").Value = 3399.25
$ws.Range("L22").Value = 4171.4116
$ws.Range("M22").Value = -3104.25
$ws.Range("N22").Value = -4761.4116

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 4024.3333
$ws.Range("I27").Value = 3399.25
$ws.Range("J27").Value = 4171.4116
$ws.Range("K27").Value = 3399.25
$ws.Range("L27").Value = 4171.4116
$ws.Range("M27").Value = -3292.25
$ws.Range("N27").Value = -4385.4116

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 4866.7896
$ws.Range("I93").Value = 5398.3335
$ws.Range("K93").Value = 5398.3335
$ws.Range("M93").Value = -4150.3335

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2646.465
$ws.Range("I136").Value = 2153.68
$ws.Range("K136").Value = 6461.039999999999
$ws.Range("M136").Value = -3911.039999999999

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 7796.689
$ws.Range("I132").Value = 7948.4326
$ws.Range("J132").Value = 7094.875
$ws.Range("K132").Value = 23845.2978
$ws.Range("L132").Value = 21284.625
$ws.Range("M132").Value = -21315.2978
$ws.Range("N132").Value = -26344.625
